$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from G1 to H1, set text "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Set H2 numeric value
$ws.Range("H2").Value = 0
